$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this block (796/797), pushing the
# existing rows 796:833 down to 798:835.
$ws.Rows("796:797").Insert()

# New record 1 (row 796) - Primera
$ws.Range("A796").Value = 3
$ws.Range("B796").Value = "Femacal de La Calera"
$ws.Range("C796").Value = "Coquimbo"
$ws.Range("D796").Value = 44939
$ws.Range("E796").Value = 5
$ws.Range("F796").Value = 100112006
$ws.Range("G796").Value = "Repollo"
$ws.Range("H796").Value = "Crespo record"
$ws.Range("I796").Value = "Primera"
$ws.Range("J796").Value = 1100
$ws.Range("K796").Value = 1200
$ws.Range("L796").Value = 1200
$ws.Range("M796").Value = 1200
$ws.Range("N796").Value = "`$/unidad"
$ws.Range("O796").Value = "Provincia de Quillota"
$ws.Range("P796").Value = 1200
$ws.Range("Q796").Value = 1
$ws.Range("R796").Value = "Hortaliza"

# New record 2 (row 797) - Segunda
$ws.Range("A797").Value = 3
$ws.Range("B797").Value = "Femacal de La Calera"
$ws.Range("C797").Value = "Coquimbo"
$ws.Range("D797").Value = 44939
$ws.Range("E797").Value = 5
$ws.Range("F797").Value = 100112006
$ws.Range("G797").Value = "Repollo"
$ws.Range("H797").Value = "Crespo record"
$ws.Range("I797").Value = "Segunda"
$ws.Range("J797").Value = 850
$ws.Range("K797").Value = 900
$ws.Range("L797").Value = 900
$ws.Range("M797").Value = 900
$ws.Range("N797").Value = "`$/unidad"
$ws.Range("O797").Value = "Provincia de Quillota"
$ws.Range("P797").Value = 900
$ws.Range("Q797").Value = 1
$ws.Range("R797").Value = "Hortaliza"
